$d = $word.ActiveDocument
$word.ActiveWindow.View.ShowHiddenBookmarks = $true
Write-Host "Bookmarks count:" $d.Bookmarks.Count
foreach ($bm in $d.Bookmarks) {
    Write-Host "Bookmark:" $bm.Name "Start:" $bm.Start "End:" $bm.End
}
Write-Host "Paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host "Para $i text: " $p.Range.Text
}
